$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.645.85'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '1.565.25'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.511'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.50%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  +5.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.245'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0898'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '1.789.68'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '1.566.63'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').Value = '28.682.40'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '227.67'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  +2.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('E30').Value = '  -3.74%  '
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = '1.406.16'
$ws.Range('E33').Value = '  +1.37%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('E36').Value = '  -1.93%  '
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.518'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.94'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.766'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('E46').Value = '  -2.05%  '
$ws.Range('D47').Value = '1.700.90'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('E48').Value = '  -6.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '84.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('E50').Value = '  +5.42%  '
$ws.Range('E51').Value = '  -0.40%  '
